# Update for release to deploy 0.1.1
# - Bump Version to 0.1.1
# - Bump Date
# - Insert a new "Jurisdiction" metadata row (blank value) after "Contact"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 0.1.0 -> 0.1.1
$ws.Range("B3").Value = "0.1.1"

# Date
$ws.Range("B8").Value = "2024-11-11T17:53:38-06:00"

# Insert a new row at 11 (pushes Description/Purpose/Copyright/... down by one)
$ws.Rows.Item(11).Insert()

# Copy formatting from the row below (now holding the old row-11 content)
# onto the freshly inserted row so it matches the table's styling.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

# Populate the new row
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
